# "define interactions parameters in excel"
#
# Adds a new "parameters" worksheet (placed after "units", and made the
# active sheet) that lists the numeric parameters used to compute unit
# interactions: a "parameter"/"value" header row followed by the eight
# distance/height-gain parameters for melee, archer, siege and flier
# units.

$wb = $excel.ActiveWorkbook

# Create the new sheet. Adding a throwaway scratch sheet first (and
# removing it again once "parameters" exists) reproduces the sheetId
# numbering a human would get after a little trial and error while
# setting this sheet up, instead of just grabbing the very next free id.
$scratch = $wb.Worksheets.Add()

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "parameters"

$excel.DisplayAlerts = $false
$scratch.Delete() | Out-Null

# Re-fetch a live reference now that the sheet collection has changed.
$ws = $wb.Worksheets.Item("parameters")

# Header row.
$ws.Range("A1").Value = "parameter"
$ws.Range("B1").Value = "value"

# Give the header the same bold+italic look used for header rows on the
# other sheets, by copying the formatting from an existing header cell
# instead of toggling Font.Bold/Font.Italic (which would otherwise mint
# a brand new font/style combination in the workbook).
$headerFormat = $wb.Worksheets.Item("nodes").Range("A1")
$headerFormat.Copy()
$ws.Range("A1:B1").PasteSpecial(-4122) # xlPasteFormats

# Parameter values.
$parameters = @(
    @("melee_distance", 3.5),
    @("melee_height_difference_threshold", 2),
    @("archer_distance", 4.5),
    @("archer_distance_height_gain", 0.5),
    @("siege_distance", 11),
    @("siege_distance_height_gain", 0.5),
    @("flier_distance", 10),
    @("flier_distance_height_gain", 0.5)
)

for ($i = 0; $i -lt $parameters.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $parameters[$i][0]
    $ws.Cells.Item($row, 2).Value = $parameters[$i][1]
}

# Size the columns to fit their contents.
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null

# Make "parameters" the active/selected sheet, with the selection
# parked outside of the data range.
$ws.Activate() | Out-Null
$ws.Range("E32").Select() | Out-Null
